$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value previously stored in C7 (now blank, like C2:C6)
$ws.Range("C7").Value = ""

# Append a new row of results from the latest script run
# Keep the date as literal text (matching A2:A7), not an auto-converted date serial
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025-03-10"
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C8").Value = "NA"
$ws.Range("D8").Value = 1
